$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fitness" (column C) values for rows 2-151 (Generation 0-149)
# according to the new run data. Rows 152-252 are unchanged.
$ws.Range("C2:C8").Value = 9729
$ws.Range("C9:C10").Value = 9673
$ws.Range("C11:C14").Value = 9660
$ws.Range("C15:C19").Value = 9004
$ws.Range("C20:C22").Value = 8062
$ws.Range("C23:C74").Value = 7639
$ws.Range("C75:C100").Value = 7581
$ws.Range("C101:C151").Value = 7295
